$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "value" header (B1) to "first_release_value".
$ws.Range("B1").Value = "first_release_value"

# Expand the single sample row into the full year-over-year series
# (yearly, 2005-12-31 .. 2025-12-31). The first and last rows only carry
# a date value (their "value" column is blank), matching the source data.
$data = @(
    @(38717, $null),
    @(39082, 12.85455285386146),
    @(39447, 8.323723251380377),
    @(39813, 2.675214973655216),
    @(40178, -14.07263945793084),
    @(40543, 14.00662378688902),
    @(40908, 10.01426242069761),
    @(41274, 5.166028195387984),
    @(41639, 0.3163778774614823),
    @(42004, 4.811464743291949),
    @(42369, 4.153025533745458),
    @(42735, 2.211964547984113),
    @(43100, 5.059349743581909),
    @(43465, 2.287635922746656),
    @(43830, 0.6462611928503614),
    @(44196, -12.47081270006417),
    @(44561, 5.490291529373104),
    @(44926, 0.01028356335206482),
    @(45291, -0.7275558254695946),
    @(45657, 0.2831497518338555),
    @(46022, $null)
)

# Carry the existing row-2 formatting (date number format, border, font,
# alignment) down to every newly inserted data row before filling values.
$lastRow = 1 + $data.Count
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A3:B" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$row = 2
foreach ($pair in $data) {
    $dateSerial = $pair[0]
    $value = $pair[1]

    $ws.Cells.Item($row, 1).Value = $dateSerial

    if ($null -ne $value) {
        $ws.Cells.Item($row, 2).Value = $value
    } else {
        $ws.Cells.Item($row, 2).ClearContents()
    }

    $row++
}
